# Daily Collection - January 2021
# Commit: "data updated till 22 Jan 2AM"
# Adds the day's (21-Jan, column AA) collection figures across the rows that
# received money that day, fills in two previously-blank agent names, and
# removes the now-stale comment that was sitting on AA56 (its cash/digital
# breakdown got superseded by the real entry).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stale comment before writing the real value into AA56.
$ws.Range("AA56").Comment.Delete()

# Fill in the two agent names that were missing.
$ws.Range("D90").Value = "Arwal"
$ws.Range("D91").Value = "PARASI"

# New collection entries for 21-Jan-2021 (column AA).
$ws.Range("AA3").Value = 3000
$ws.Range("AA9").Value = 2000
$ws.Range("AA10").Value = 5000
$ws.Range("AA14").Value = 1000
$ws.Range("AA22").Value = 1000
$ws.Range("AA26").Value = 1050
$ws.Range("AA31").Value = 1000
$ws.Range("AA34").Value = 1000
$ws.Range("AA41").Value = 3000
$ws.Range("AA44").Value = 3000
$ws.Range("AA54").Value = 4000
$ws.Range("AA55").Value = 5000
$ws.Range("AA56").Value = 2000
$ws.Range("AA59").Value = 2000
$ws.Range("AA66").Value = 6000
$ws.Range("AA69").Value = 1000
$ws.Range("AA90").Value = 5000
$ws.Range("AA91").Value = 5000
$ws.Range("AA95").Value = 2000

# Move the frozen-pane cursor to where today's entry work left off.
$ws.Range("AA67").Select()
